{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst verNoJupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightMarker = \"Contact: luizeleno@usp.br\";\n\nlet verNoJupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === verNoJupiterText) {\n    verNoJupiterIndex = i;\n  } else if (t.indexOf(copyrightMarker) !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (verNoJupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the target paragraphs to remove.\");\n}\n\n// Also remove the blank spacer paragraph that sits right after the\n// copyright line (it collapses together with the two removed lines in\n// the canonical edit). Delete from the highest index downward so the\n// earlier indices stay valid.\nconst indexesToDelete = [verNoJupiterIndex, copyrightIndex, copyrightIndex + 1];\nindexesToDelete.sort((a, b) => b - a);\n\nfor (const idx of indexesToDelete) {\n  paragraphs.items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$verNoJupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightMarker = \"Contact: luizeleno@usp.br\"\n\n$idxVerNoJupiter = -1\n$idxCopyright = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($t -eq $verNoJupiterText) {\n        $idxVerNoJupiter = $i\n    } elseif ($t.Contains($copyrightMarker)) {\n        $idxCopyright = $i\n    }\n}\n\nif ($idxVerNoJupiter -eq -1 -or $idxCopyright -eq -1) {\n    throw \"Could not locate the target paragraphs to remove.\"\n}\n\n# Remove, highest index first so earlier indices stay valid:\n#  - the blank spacer paragraph right after the copyright line\n#  - the copyright line itself\n#  - the \"Ver no Jupiter...\" line\n$d.Paragraphs.Item($idxCopyright + 1).Range.Delete()\n$d.Paragraphs.Item($idxCopyright).Range.Delete()\n$d.Paragraphs.Item($idxVerNoJupiter).Range.Delete()\n"}
